$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "20.474.93"
$ws.Range("E2").Value = "  +2.39%  "

# Row 3
$ws.Range("D3").Value = "1.464.29"
$ws.Range("E3").Value = "  +3.27%  "

# Row 4
$ws.Range("D4").Value = "1.011"
$ws.Range("E4").Value = "  +0.96%  "

# Row 5
$ws.Range("D5").Value = "0.9248"
$ws.Range("E5").Value = "  -7.61%  "

# Row 6
$ws.Range("D6").Value = "280.62"
$ws.Range("E6").Value = "  +2.44%  "

# Row 7
$ws.Range("D7").Value = "0.3711"
$ws.Range("E7").Value = "  +0.33%  "

# Row 8
$ws.Range("D8").Value = "0.3183"
$ws.Range("E8").Value = "  +3.25%  "

# Row 9
$ws.Range("D9").Value = "40.78"
$ws.Range("E9").Value = "  +3.08%  "

# Row 10
$ws.Range("D10").Value = "1.052"
$ws.Range("E10").Value = "  +4.74%  "

# Row 11
$ws.Range("D11").Value = "0.06639"
$ws.Range("E11").Value = "  +0.75%  "

# Row 12
$ws.Range("D12").Value = "1.006"

# Row 13
$ws.Range("D13").Value = "5.563"
$ws.Range("E13").Value = "  +2.00%  "

# Row 14
$ws.Range("D14").Value = "18.09"
$ws.Range("E14").Value = "  +5.97%  "

# Row 15
$ws.Range("D15").Value = "6.230"
$ws.Range("E15").Value = "  +0.86%  "

# Row 16
$ws.Range("D16").Value = "1.474.61"
$ws.Range("E16").Value = "  +3.39%  "

# Row 17
$ws.Range("D17").Value = "0.00001035"
$ws.Range("E17").Value = "  +2.64%  "

# Row 18
$ws.Range("D18").Value = "0.9249"
$ws.Range("E18").Value = "  -7.64%  "

# Row 19
$ws.Range("D19").Value = "0.05717"
$ws.Range("E19").Value = "  -0.69%  "

# Row 20
$ws.Range("D20").Value = "71.70"
$ws.Range("E20").Value = "  -3.43%  "

# Row 21
$ws.Range("D21").Value = "5.689"
$ws.Range("E21").Value = "  +1.04%  "

# Row 22
$ws.Range("D22").Value = "14.70"
$ws.Range("E22").Value = "  +1.42%  "

# Row 23
$ws.Range("D23").Value = "11.19"
$ws.Range("E23").Value = "  +2.08%  "

# Row 24
$ws.Range("D24").Value = "2.292"
$ws.Range("E24").Value = "  -1.73%  "

# Row 25
$ws.Range("D25").Value = "20.581.87"
$ws.Range("E25").Value = "  +2.88%  "

# Row 26
$ws.Range("D26").Value = "2.302"
$ws.Range("E26").Value = "  +1.07%  "

# Row 27
$ws.Range("D27").Value = "137.86"
$ws.Range("E27").Value = "  -1.26%  "

# Row 28
$ws.Range("D28").Value = "17.48"
$ws.Range("E28").Value = "  +2.87%  "

# Row 29
$ws.Range("D29").Value = "1.636.66"
$ws.Range("E29").Value = "  +3.44%  "

# Row 30
$ws.Range("D30").Value = "113.57"
$ws.Range("E30").Value = "  +3.98%  "

# Row 31
$ws.Range("D31").Value = "3.950"
$ws.Range("E31").Value = "  +1.47%  "

# Row 32
$ws.Range("D32").Value = "5.271"
$ws.Range("E32").Value = "  -2.45%  "

# Row 33
$ws.Range("D33").Value = "0.8466"
$ws.Range("E33").Value = "  -2.09%  "

# Row 34
$ws.Range("B34").Value = "WEMIXTOKEN"
$ws.Range("C34").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D34").Value = "1.578"
$ws.Range("E34").Value = "  +23.44%  "

# Row 35
$ws.Range("B35").Value = "Stellar"
$ws.Range("C35").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D35").Value = "0.07804"
$ws.Range("E35").Value = "  +0.94%  "

# Row 36
$ws.Range("D36").Value = "0.06082"
$ws.Range("E36").Value = "  +5.82%  "

# Row 37
$ws.Range("D37").Value = "4.879"
$ws.Range("E37").Value = "  +1.96%  "

# Row 38
$ws.Range("B38").Value = "TrustWalletToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D38").Value = "1.139"
$ws.Range("E38").Value = "  +6.33%  "

# Row 39
$ws.Range("B39").Value = "Aptos"
$ws.Range("C39").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D39").Value = "10.64"
$ws.Range("E39").Value = "  -0.73%  "

# Row 40
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").Value = "0.02065"
$ws.Range("E40").Value = "  +1.34%  "

# Row 41
$ws.Range("D41").Value = "0.1889"
$ws.Range("E41").Value = "  -2.07%  "

# Row 42
$ws.Range("D42").Value = "0.9436"
$ws.Range("E42").Value = "  -5.78%  "

# Row 43
$ws.Range("D43").Value = "7.474"
$ws.Range("E43").Value = "  -11.60%  "

# Row 44
$ws.Range("D44").Value = "0.5391"
$ws.Range("E44").Value = "  +1.54%  "

# Row 45
$ws.Range("B45").Value = "PancakeSwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D45").Value = "3.585"
$ws.Range("E45").Value = "  +1.45%  "

# Row 46
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "12.46"
$ws.Range("E46").Value = "  +1.72%  "

# Row 47
$ws.Range("D47").Value = "123.40"
$ws.Range("E47").Value = "  +12.55%  "

# Row 48
$ws.Range("D48").Value = "0.5300"
$ws.Range("E48").Value = "  +3.24%  "

# Row 49
$ws.Range("D49").Value = "1.827"
$ws.Range("E49").Value = "  +1.16%  "

# Row 50
$ws.Range("D50").Value = "0.06436"
$ws.Range("E50").Value = "  +4.56%  "

# Row 51
$ws.Range("D51").Value = "1.042"
$ws.Range("E51").Value = "  -0.85%  "
